$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (PortraitPath), shifting PortraitPath and PathFile right.
$ws.Columns("H").Insert()

$ws.Range("H1").Value = "SlotCount"
$ws.Range("H2").Value = "int32"

for ($r = 3; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

[void]$ws.Range("H3:H15").Select()
